$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for 2026/02/01, appended after the last existing row (82).
$newRow = 83

# Copy formatting from the last data row so the new row matches the
# existing table's style (centered alignment, same column widths, etc.).
$ws.Range("A82:C82").Copy()
$ws.Range("A83:C83").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column A holds a date written as plain text ("2026/02/01"), not a real
# date value. Assigning that text via .Formula (as a quoted string literal)
# keeps Excel's automatic date-recognition from turning it into a date
# serial number, then PasteSpecial-Values collapses the formula down to
# its literal text result while preserving the row's number format/style.
$ws.Cells.Item($newRow, 1).Formula = "=""2026/02/01"""
$ws.Range("A83").Copy()
$ws.Range("A83").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1168
